$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

for ($r = 2; $r -le 13; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null -and $val.ToString().StartsWith("REPSWITCH1_Practice/")) {
        $cell.Value2 = $val.ToString().Replace("REPSWITCH1_Practice/", "Pictures_Practice/")
    }
}
